$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 377; this shifts rows 377-400 down to 378-401
# and extends the used range to A1:T401, matching the target diff.
$ws.Rows.Item(377).Insert()

# Populate the newly inserted row 377 with the new weekly record.
$ws.Cells.Item(377, 1).Value2 = 4
$ws.Cells.Item(377, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(377, 3).Value2 = "Los Lagos"
$ws.Cells.Item(377, 4).Value2 = 45021
$ws.Cells.Item(377, 5).Value2 = 10
$ws.Cells.Item(377, 6).Value2 = "Fruta"
$ws.Cells.Item(377, 7).Value2 = 100108
$ws.Cells.Item(377, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(377, 9).Value2 = 100108005
$ws.Cells.Item(377, 10).Value2 = "Piña"
$ws.Cells.Item(377, 11).Value2 = "Caramelo"
$ws.Cells.Item(377, 12).Value2 = "Primera"
$ws.Cells.Item(377, 13).Value2 = 200
$ws.Cells.Item(377, 14).Value2 = 21000
$ws.Cells.Item(377, 15).Value2 = 22000
$ws.Cells.Item(377, 16).Value2 = 21500
$ws.Cells.Item(377, 17).Value2 = "$/caja 12 unidades"
$ws.Cells.Item(377, 18).Value2 = "Ecuador"
$ws.Cells.Item(377, 19).Value2 = 1792
$ws.Cells.Item(377, 20).Value2 = 12
